$d = $word.ActiveDocument

# Locate the three paragraphs to remove:
#   1) the empty paragraph right after "Artigos sobre metodologias..."
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "(c) 2020 . Contact: ..." paragraph
# and delete the whole range spanning them (start of the empty paragraph
# through the end of the copyright paragraph, including its paragraph mark).

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "Ver no Jupiter Salvar em pdf Salvar em docx`r") {
        $startPara = $d.Paragraphs.Item($i - 1)
        $endPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

$delStart = $startPara.Range.Start
$delEnd = $endPara.Range.End

$r = $d.Range($delStart, $delEnd)
$r.Delete()
